$d = $word.ActiveDocument

# The edit rewrites the three paragraphs of the body:
#  1) "Anexo 1" (Titulo1)            -> "Appendix A: Some clarification", split into
#     several runs, with a new bookmark (_Hlk153205755) wrapped around "Appendix ",
#     and the existing "introducción" bookmarkStart kept in place.
#  2) "Quis autem vel eum iure reprehenderit " (FirstParagraph) -> proofErr marks
#     removed and the runs merged into two runs; bookmarkEnd id=0 kept at the end.
#  3) The empty "Textoindependiente" paragraph is kept as-is.
# In all three paragraphs the stray <w:rPr><w:lang w:val="es-ES"/></w:rPr> that
# Word had attached to the paragraph mark (pPr) and to every run is dropped.

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Ttulo1"/>
            </w:pPr>
            <w:bookmarkStart w:id="0" w:name="introducción"/>
            <w:bookmarkStart w:id="1" w:name="_Hlk153205755"/>
            <w:r><w:t>A</w:t></w:r>
            <w:r><w:t>ppendix</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:bookmarkEnd w:id="1"/>
            <w:r><w:t xml:space="preserve">A: </w:t></w:r>
            <w:r><w:t>Some clarification</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="FirstParagraph"/>
            </w:pPr>
            <w:r><w:t>Quis autem vel eum iure reprehenderit</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Textoindependiente"/>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($xml)
